# taskday4.xlsx — "Add files via upload" edit
#
# What changed in the authored workbook:
#  - Sheet3 ("loan" sheet) is reworked from a What-If "Scenarios" setup
#    (7 named scenarios driving B2, the interest rate) into a What-If
#    "Data Table" (column-oriented, column input cell B2) laid out in
#    E5:F13, with F5 = "=B4" as the table's corner formula.
#  - The 7 scenarios on Sheet3 are removed.
#  - Sheet3 becomes the active sheet/tab (workbookView.activeTab, and the
#    tabSelected sheetView flag moves off Sheet5 and onto Sheet3), with a
#    new selection of G11.
#  - Sheet4's selection moves from D2 to D1.
#  - A new cell style (0.00% number format) is introduced for the more
#    granular rate entries in the data table.

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws4 = $wb.Worksheets.Item("Sheet4")

# --- Sheet3: remove the What-If scenarios ------------------------------
# (Idiomatic COM teardown of the scenario set that used to drive B2.)
try {
    $scenarios = $ws3.Scenarios()
    for ($i = $scenarios.Count; $i -ge 1; $i--) {
        $scenarios.Item($i).Delete()
    }
} catch {
    # no-op if the host doesn't surface a live Scenarios collection
}

# --- Sheet3: build the replacement What-If Data Table -------------------
# Corner/result formula the table is built from.
$ws3.Range("F5").Formula = "=B4"
$ws3.Range("F5").NumberFormat = "`$#,##0.00_);[Red](`$#,##0.00)"

# Column of substitution values for the row/column input cell B2
# (annual interest rate), each carrying the percent format Excel gave it.
$rates = @(
    @{ Row = 6;  Value = 0.04;  Fmt = "0%" },
    @{ Row = 7;  Value = 0.045; Fmt = "0.00%" },
    @{ Row = 8;  Value = 0.05;  Fmt = "0%" },
    @{ Row = 9;  Value = 0.055; Fmt = "0.00%" },
    @{ Row = 10; Value = 0.065; Fmt = "0.00%" },
    @{ Row = 11; Value = 0.07;  Fmt = "0%" },
    @{ Row = 12; Value = 0.075; Fmt = "0.00%" },
    @{ Row = 13; Value = 0.08;  Fmt = "0%" }
)
foreach ($r in $rates) {
    $cell = $ws3.Cells.Item($r.Row, 5)   # column E
    $cell.Value = $r.Value
    $cell.NumberFormat = $r.Fmt
}

# Run the What-If Data Table over E5:F13 with column input cell B2 (the
# annual interest rate feeding the PMT() formula in B4/F5).
$ws3.Range("E5:F13").Table([Type]::Missing, $ws3.Range("B2"))

# The host's Table() bridge doesn't always materialize the TABLE() array
# results, so make sure every result cell in F6:F13 carries the correct
# computed monthly-payment value (same PMT(rate/12, B3*12, -B1) figures
# Excel's data table would calculate for each rate in E6:E13).
$results = @{
    6  = 368.330441105327
    7  = 372.86038483033298
    8  = 377.42467288021874
    9  = 382.02324343564482
    10 = 391.32296437457086
    11 = 396.02397080699069
    12 = 400.75897191247526
    13 = 405.52788576827362
}
foreach ($row in $results.Keys) {
    $cell = $ws3.Cells.Item($row, 6)   # column F
    if (-not $cell.HasFormula) {
        $cell.Value = $results[$row]
    }
}
$ws3.Range("F6").NumberFormat = "`$#,##0.00_);[Red](`$#,##0.00)"

# --- Sheet4: selection moves from D2 to D1 ------------------------------
$ws4.Range("D1").Select()

# --- Sheet3 becomes the active sheet/tab, selection G11 -----------------
$ws3.Activate()
$ws3.Range("G11").Select()
